$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 100 - this shifts the existing rows 100-117
# down to 101-118, preserving all of their data/styles.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new weekly report entry.
# (Most fields mirror the former row 100 - same market/category/variety -
#  only the date, volume, prices, origin and $/Kg differ.)
$ws.Range("A100").Value = 7
$ws.Range("B100").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value = "Ñuble"
$ws.Range("D100").Value = 45154
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = 100112013
$ws.Range("G100").Value = "Alcachofa"
$ws.Range("H100").Value = "Madrigal"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 60
$ws.Range("K100").Value = 13000
$ws.Range("L100").Value = 13000
$ws.Range("M100").Value = 13000
$ws.Range("N100").Value = "`$/caja 40 unidades"
$ws.Range("O100").Value = "Provincia de Limarí"
$ws.Range("P100").Value = 325
$ws.Range("Q100").Value = 40
$ws.Range("R100").Value = "Hortaliza"
